# Apply the edits described by the OOXML diff to Hoja1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Colors used as literal RGB fills that already exist in the style table ---
$green  = 5296274   # RGB(0x92,0xD0,0x50) -> fill used by styles 3 and 7
$yellow = 65535      # RGB(0xFF,0xFF,0x00) -> fill used by style 4

# --- Row 9 / 11: E,F switch from the plain red style to red+green style ---
$ws.Range("E9").Interior.Color = $green
$ws.Range("F9").Interior.Color = $green
$ws.Range("E11").Interior.Color = $green
$ws.Range("F11").Interior.Color = $green

# --- Row 10: E,F gain the plain green style (no font change) ---
$ws.Range("E10").Interior.Color = $green
$ws.Range("F10").Interior.Color = $green

# --- Row 12: turn the blank banner row into the header-label row (like row 1) ---
$ws.Range("A12:H12").Font.Bold = $false
$ws.Range("A12:H12").Interior.Color = $yellow
$ws.Range("A12").Value = "bunny.pcd"
$ws.Range("B12").Value = "/10"
$ws.Range("C12").Value = "/5"
$ws.Range("D12").Value = "/2"
$ws.Range("E12").Value = "default"
$ws.Range("F12").Value = "x2"
$ws.Range("G12").Value = "x5"
$ws.Range("H12").Value = "x10"

# --- Row 20: B gains the plain green style ---
$ws.Range("B20").Interior.Color = $green

# --- Row 21 / 22: B gains the red+green style ---
$ws.Range("B21").Interior.Color = $green
$ws.Range("B22").Interior.Color = $green

# --- Row 23: blank banner row -> yellow style, with A/E labelled ---
$ws.Range("A23:H23").Font.Bold = $false
$ws.Range("A23:H23").Interior.Color = $yellow
$ws.Range("A23").Value = "bunny.pcd"
$ws.Range("E23").Value = "default"

# --- Row 31: C gains the red+green style (font was plain, needs to become red too) ---
$ws.Range("C31").Font.Color = 255
$ws.Range("C31").Interior.Color = $green

# --- Row 32 / 33: C gains the red+green style ---
$ws.Range("C32").Interior.Color = $green
$ws.Range("C33").Interior.Color = $green

# --- Row 34: blank banner row -> yellow style, with A/E labelled, extend to I/J ---
$ws.Range("A34:H34").Font.Bold = $false
$ws.Range("A34:H34").Interior.Color = $yellow
$ws.Range("A34").Value = "bunny.pcd"
$ws.Range("E34").Value = "default"
$ws.Range("I34").Interior.Color = $yellow
$ws.Range("J34").Interior.Color = $yellow

# --- Row 42: G,H gain the plain green style ---
$ws.Range("G42").Interior.Color = $green
$ws.Range("H42").Interior.Color = $green

# --- Row 43 / 44: G,H gain the red+green style ---
$ws.Range("G43").Interior.Color = $green
$ws.Range("H43").Interior.Color = $green
$ws.Range("G44").Interior.Color = $green
$ws.Range("H44").Interior.Color = $green

# --- Row 45: blank banner row -> full header-label row (like row 1/12), shifted one column ---
$ws.Range("A45:J45").Font.Bold = $false
$ws.Range("A45:J45").Interior.Color = $yellow
$ws.Range("A45").Value = "bunny.pcd"
$ws.Range("B45").Value = "x0"
$ws.Range("C45").Value = "/10"
$ws.Range("D45").Value = "/5"
$ws.Range("E45").Value = "/2"
$ws.Range("F45").Value = "default"
$ws.Range("G45").Value = "x2"
$ws.Range("H45").Value = "x5"
$ws.Range("I45").Value = "x10"

# --- Sheet view: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("B45").Select()
